$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New, more realistic/complete demo data for the variable_modality table.
$data = @(
    @("ser_pub_loc___variable_3",   "modality_1"),
    @("ser_pub_loc___variable_5",   "modality_1"),
    @("accident_route___variable_2","a_or_b"),
    @("ser_pub_loc___canton",       "canton_sigle"),
    @("ser_pub_loc___langue",       "langue_sigle"),
    @("ser_pub_loc___nouveau",      "oui_non"),
    @("ser_pub_loc___nouveau",      "vide")
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Grow the Excel Table (ListObject) so it covers the newly added rows.
$table = $ws.ListObjects.Item(1)
$table.Resize($ws.Range("A1:B8"))

# Select B9 (the cell right below/after the new last table row) like in the
# final authored workbook.
$ws.Range("B9").Select()
